$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Classic Multiplication Time" (row 2) and
# "Strassen Multiplication Time" (row 3) rows - only a single
# timing row remains after the fix.
$ws.Rows("2:3").Delete()

# Update the remaining row's label and timing value.
$ws.Range("A1").Value = "Classic Multiplication Time"
$ws.Range("B1").Value = 0.001169
